$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Price/Volume cell in this sheet is stored as text (values such as
# "26.060.37" or "  -0.07%  " are not valid Excel numbers, and plain-numeric
# -looking ones like "218.10" or "5.900" must keep their exact, sometimes
# trailing-zero, string form). Assigning a numeric-looking string straight
# to .Value would let Excel auto-convert it to a real number (dropping
# trailing zeros, e.g. "5.900" -> 5.9), so every touched cell is first
# forced to Text format, written, then restored to the workbook's default
# (Normal) style so no formatting differences are introduced.
$touchedCells = "D2", "D3", "D5", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51", "E2", "E3", "E4", "E5", "E6", "E8", "E9", "E10", "E11", "E12", "E13", "E14", "E15", "E16", "E18", "E19", "E20", "E21", "E22", "E23", "E24", "E25", "E26", "E27", "E28", "E29", "E30", "E31", "E32", "E33", "E34", "E35", "E36", "E37", "E38", "E39", "E40", "E42", "E43", "E44", "E45", "E46", "E48", "E50", "E51"
foreach ($addr in $touchedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.061.24"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.650.73"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "218.13"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D8").Value = "0.2635"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").Value = "0.06327"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "20.35"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "0.07666"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "4.583"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").Value = "1.634.72"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "1.877.32"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "0.5594"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "0.0₅8136"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "65.22"
$ws.Range("D18").Value = "26.048.78"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "10.49"
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("D22").Value = "191.45"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "5.900"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "143.96"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "0.1187"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").Value = "7.200"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "15.87"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  +2.74%  "
$ws.Range("D30").Value = "0.05432"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("D31").Value = "1.266"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "3.444"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "3.342"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "1.554"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").Value = "2.423"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("D36").Value = "2.780"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").Value = "0.9446"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "0.5631"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "0.01579"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "5.864"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D42").Value = "1.026.42"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("D43").Value = "0.8268"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "100.85"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("D45").Value = "1.785.05"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("D47").Value = "57.30"
$ws.Range("D48").Value = "0.9982"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "0.4330"
$ws.Range("D50").Value = "7.927"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "0.05141"
$ws.Range("E51").Value = "  -3.48%  "

foreach ($addr in $touchedCells) {
    $ws.Range($addr).Style = "Normal"
}
